$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 with the new person's data (strings entered first so shared
# string table is built in the expected order)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "MEJIA ARANGO"
$ws.Range("C3").Value = "ISABELLA MARIA"
$ws.Range("D3").Value = "T1019906212"
$ws.Range("E3").Value = 1

# Update row 2: ID_STATUS 3 -> 1, STATUS "Inactivo(a)" -> "Activo(a)"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "Activo(a)"

$ws.Range("F3").Value = "Activo(a)"
$ws.Range("G3").Value = 1019906212
$ws.Range("H3").Value = "F Mejia"
$ws.Range("I3").Value = "1019906212.jpg"

# Set column I width to match diff (stored width 16, customWidth)
$ws.Columns.Item(9).ColumnWidth = 15.17

# Update selection to match diff (active cell F2)
$ws.Range("F2").Select()
